# "Add files via upload" - updated event descriptions in the
# "Análise dos Eventos para cada Cenário" table to reflect the vendor's
# (not the client's) perspective of the use case.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = "Efetuar a venda do produto"
$ws.Range("E4").Value = "Receber o pagamento em dinheiro"
$ws.Range("E5").Value = "Receber o pagamento em cartão de crédito/debito"
$ws.Range("E6").Value = "Tratar a resposta da operação realizada"

# Column E needs to fit the new (longer) text.
$ws.Columns.Item(5).EntireColumn.AutoFit() | Out-Null

# Reflect the last active cell/selection recorded in the saved file.
$ws.Range("G3").Select()
